$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Fill in the speaker notes (currently empty) for slides 1-12.
# -----------------------------------------------------------------
$notes = @{
    1  = "Today I did my presentation on Page Ranking. This presentation covers the definition of page ranking, its origins, the TF-IDF calculations, and a couple of examples. This presentation will also briefly go over different kinds of page ranking methods."
    2  = "First off, what exactly is page ranking and why is it so important? Let us start by considering what a search query is. By now, everyone has made a search on a search engine such as Google. The search engine will try to return the results of the highest quality, which often translates into what is the most viewed link there. PageRank, the term, was actually invented by Google’s founders, Larry Page and Sergey Brin. This particular ranking system evaluates the quality and quantity of links to a page, determining overall relativity on a scale of 1 to 10."
    3  = "Before we talk about PageRank, we should take a look at two different terms as well as the combination of the two which had a part in determining page rankings. The first is TF, which is known as term frequency. This pertains to how often a term appears in a document. The calculation is the number of times a term t is present in a document divided by the total number of terms in the document. The second term is IDF, inverse document frequency, measuring how important a term is, weighing down the frequency of terms while scaling up the rare ones. The calculation is the natural log of the total number of documents divided by the number of documents with the term t in it."
    4  = "Here is an example problem. The page I referenced these terms from had a similar problem to this. A document has 500 words, and the word faucet appears in the document 12 times. 12/500 = 0.024. This is the TF. Now say that there are 1000 documents, with the term faucet appearing in 150 of these documents. The IDF calculation is log(e)(1000/150) which is approximately equal to 0.824."
    5  = "The combined term is simply TF-IDF. Despite what it looks like, this is not subtracting the IDF from the TF. This term is simply the multiplication of both. Taking in the outputs of the TF and IDF calculations, the TF-IDF calculations are 0.024 * 0.824, which is approximately equal to 0.00198. I try to keep significant digits to three here. Now, TF-IDF is the important deal for page ranking, reflecting how important a term is to a document, or perhaps a collection. This helps in text mining, info retrieval, user modeling. TF-IDF can also be used for stop-word filtering for fields such as text summarization and classification"
    6  = "Now I will talk a bit more about Google PageRank, and also about Alexa, which is a more complex form of page ranking."
    7  = "As said before, this was founded by Larry Page and Sergey Brin. It works by counting the quantity and the quality of links to a page, determining the estimate of how important a website is for query results. The algorithm is pretty rough to follow so let us go over soome of the terms. PR(A) is the PageRank of a page A, and PR(Ti) is the PageRank of pages Ti, all of which lead to page A. C(Ti) equal the number of outbound links on page Ti, and d is actually a damping factor that only accepts two values, 0 and 1"
    8  = "For an example to this, I will use myspace.com. After inputing it in hte text box shown in the image and clicking Check PR, I find that the score between 0 and 10 is 8. The conclusion I found from this experiment is this; the more outbound links a page T will have, the less that a page A will benefit from page T. Weighted PageRanks are then added up, then the damping factor will come into play."
    9  = "Alexa is a page ranking tool which is more advanced than Google’s PR, getting in more than just the page rank. Owned by Amazon, Alexa will get things such as website traffic, statistics, popularity, visitor metrics, demographics, geography, upstream sites, linking sites, related sites, and even how a site loads."
    10 = "So again, I use Myspace.com as an example. There is a traffic ranking towards other sites, a global rank for everyone, and a rank for a certain country where the user got the ranking from. Visitor metrics, basically they deal with how often a site gets visited monthly. Audience geography shows the location by country for those who visit a site, so expect traffic to be frequent in places such as the United States and China. Upstream sites deal with the sites that are used to reach a site, which will almost always have Google at the top. Alexa also asks where do site visitors go next, as well as other sites that link to a site. Not to forget related sites. Surprisingly, Facebook did not show up next to Myspace. Demographics show overall gender, browsing location, and education of viewers."
    11 = "Well, to summarize page ranking is a major ranking system for Google and other engines to determine relevancy and popularity of websites given from a search query. Alexa and Google PR are simply two of the tools that can be used. Other ones do exist, but I go over two of the more known ones. Be mindful there are multiple ways of determining a page’s rank. Finally, TF and IDF are important to know just how these page ranking tools work, and are a crude way of determining a page rank via a program."
    12 = "Questions?"
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $notesShape = $slide.NotesPage.Shapes.Item(2)
    $notesShape.TextFrame.TextRange.Text = $notes[$i]
}

# -----------------------------------------------------------------
# 2) Slide 9 ("Alexa"): append a hyperlink line + trailing blank
#    paragraph to the second body placeholder, and tighten the
#    spacing after the first paragraph.
# -----------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$alexaShape = $slide9.Shapes.Item(3)
$alexaRange = $alexaShape.TextFrame.TextRange

$existingText = "Measures popularity, visitor metrics, audience geography, upstream sites, linking sites, related sites, how fast a site loads, audience demographics"
$linkText = "https://www.alexa.com/siteinfo"

$alexaRange.Text = $existingText + "`r" + $linkText + "`r"

$para1 = $alexaRange.Paragraphs(1, 1)
$para2 = $alexaRange.Paragraphs(2, 1)
$para3 = $alexaRange.Paragraphs(3, 1)

$para1.ParagraphFormat.SpaceBefore = 0
$para1.ParagraphFormat.SpaceAfter = 0

$para2.ParagraphFormat.SpaceBefore = 16
$para2.ParagraphFormat.SpaceAfter = 0
$para2.Font.Underline = $true
$para2.Font.Color.ObjectThemeColor = 11
$linkSettings = $para2.ActionSettings(1)
$linkSettings.Hyperlink.Address = $linkText

$para3.ParagraphFormat.SpaceBefore = 16
$para3.ParagraphFormat.SpaceAfter = 16
